$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.06, 0.04000000000000001, 8.919944286346436, 3.40625),
    @(39.78, 0.16, 24.38289165496826, 15.984375),
    @(39.78, 0.16, 25.35724759101868, 13.09375),
    @(39.06, 0.04000000000000001, 6.70301628112793, 3.40625),
    @(39.06, 0.04000000000000001, 7.414100885391235, 2.921875)
)

$startRow = 215
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
